$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: I1 = "I0", J1 = "IF", styled like the other headers (bold/border/center)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Data values for I2:I66 and J2:J66
$iValues = @(3,6,9,8,8,7,6,5,7,6,7,7,8,7,8,9,7,7,7,7,10,8,8,8,7,6,6,5,6,8,6,7,7,8,6,6,5,7,8,9,7,8,6,8,7,9,8,8,6,7,4,6,6,7,6,9,7,5,8,9,4,4,5,5,3)
$jValues = @(4,6,9,8,8,7,6,5,8,7,7,7,8,7,8,9,7,8,7,7,10,8,8,8,7,7,6,6,7,8,6,7,7,8,7,6,5,7,8,9,7,8,6,8,7,9,9,8,7,7,4,6,6,8,7,9,7,6,8,9,4,4,5,5,3)

for ($r = 2; $r -le 66; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
